$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Volume Number text (A8): "...Number  46" -> "...Number  47" ---
$rA8 = $ws.Range("A8")
$textA8 = $rA8.Characters().Text
$idxA8 = $textA8.LastIndexOf("46")
$rA8.Characters($idxA8 + 1, 2).Text = "47"

# --- Update report week date range (C9): 11/13/2023-11/19/2023 -> 11/20/2023-11/26/2023 ---
$rC9 = $ws.Range("C9")
$textC9 = $rC9.Characters().Text
$idxC9a = $textC9.IndexOf("11/13/2023")
$rC9.Characters($idxC9a + 1, 10).Text = "11/20/2023"
$textC9b = $rC9.Characters().Text
$idxC9b = $textC9b.IndexOf("11/19/2023")
$rC9.Characters($idxC9b + 1, 10).Text = "11/26/2023"

# --- Update weekly crime-stat data table (rows 14-30) with newly collected figures ---
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = -16.666666666666
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 116
$ws.Range("J14").Value = 118
$ws.Range("K14").Value = -1.694915254237
$ws.Range("L14").Value = -12.781954887218
$ws.Range("M14").Value = -1.694915254237
$ws.Range("N14").Value = -74.837310195227
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 75
$ws.Range("F15").Value = 17
$ws.Range("G15").Value = 21
$ws.Range("H15").Value = -19.047619047619
$ws.Range("I15").Value = 330
$ws.Range("J15").Value = 350
$ws.Range("K15").Value = -5.714285714285
$ws.Range("L15").Value = -0.302114803625
$ws.Range("M15").Value = 21.771217712177
$ws.Range("N15").Value = -49.695121951219
$ws.Range("C16").Value = 103
$ws.Range("D16").Value = 96
$ws.Range("E16").Value = 7.291666666666
$ws.Range("F16").Value = 368
$ws.Range("G16").Value = 436
$ws.Range("H16").Value = -15.596330275229
$ws.Range("I16").Value = 4439
$ws.Range("J16").Value = 4672
$ws.Range("K16").Value = -4.987157534246
$ws.Range("L16").Value = 24.551066217732
$ws.Range("M16").Value = 5.791229742612
$ws.Range("N16").Value = -70.093646836892
$ws.Range("C17").Value = 117
$ws.Range("D17").Value = 149
$ws.Range("E17").Value = -21.476510067114
$ws.Range("F17").Value = 603
$ws.Range("G17").Value = 597
$ws.Range("H17").Value = 1.005025125628
$ws.Range("I17").Value = 7388
$ws.Range("J17").Value = 6717
$ws.Range("K17").Value = 9.989578680958
$ws.Range("L17").Value = 29.319096796779
$ws.Range("M17").Value = 81.256133464180
$ws.Range("N17").Value = -10.491882723528
$ws.Range("C18").Value = 51
$ws.Range("D18").Value = 47
$ws.Range("E18").Value = 8.510638297872
$ws.Range("F18").Value = 199
$ws.Range("G18").Value = 202
$ws.Range("H18").Value = -1.485148514851
$ws.Range("I18").Value = 2630
$ws.Range("J18").Value = 2665
$ws.Range("K18").Value = -1.313320825515
$ws.Range("L18").Value = 31.434282858570
$ws.Range("M18").Value = -13.429888084266
$ws.Range("N18").Value = -84.623479887745
$ws.Range("C19").Value = 134
$ws.Range("D19").Value = 142
$ws.Range("E19").Value = -5.633802816901
$ws.Range("F19").Value = 617
$ws.Range("G19").Value = 618
$ws.Range("H19").Value = -0.161812297734
$ws.Range("I19").Value = 7283
$ws.Range("J19").Value = 7339
$ws.Range("K19").Value = -0.763046736612
$ws.Range("L19").Value = 18.345791355216
$ws.Range("M19").Value = 71.243827886198
$ws.Range("N19").Value = 5.949956357288
$ws.Range("C20").Value = 66
$ws.Range("D20").Value = 75
$ws.Range("E20").Value = -12
$ws.Range("F20").Value = 329
$ws.Range("G20").Value = 290
$ws.Range("H20").Value = 13.448275862069
$ws.Range("I20").Value = 4693
$ws.Range("J20").Value = 3516
$ws.Range("K20").Value = 33.475540386803
$ws.Range("L20").Value = 72.600220669363
$ws.Range("M20").Value = 145.706806282723
$ws.Range("N20").Value = -66.442617089739
$ws.Range("C21").Value = 483
$ws.Range("D21").Value = 519
$ws.Range("E21").Value = -6.936416184971
$ws.Range("F21").Value = 2141
$ws.Range("G21").Value = 2172
$ws.Range("H21").Value = -1.427255985267
$ws.Range("I21").Value = 26879
$ws.Range("J21").Value = 25377
$ws.Range("K21").Value = 5.918745320565
$ws.Range("L21").Value = 30.385641523162
$ws.Range("M21").Value = 50.481469040421
$ws.Range("N21").Value = -56.770188333306
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 29
$ws.Range("G22").Value = 24
$ws.Range("H22").Value = 20.833333333333
$ws.Range("I22").Value = 288
$ws.Range("J22").Value = 329
$ws.Range("K22").Value = -12.462006079027
$ws.Range("L22").Value = 12.062256809338
$ws.Range("M22").Value = 3.225806451612
$ws.Range("C23").Value = 20
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = -37.5
$ws.Range("F23").Value = 112
$ws.Range("G23").Value = 128
$ws.Range("H23").Value = -12.5
$ws.Range("I23").Value = 1577
$ws.Range("J23").Value = 1447
$ws.Range("K23").Value = 8.984105044920
$ws.Range("L23").Value = 38.820422535211
$ws.Range("M23").Value = 61.247443762781
$ws.Range("C24").Value = 237
$ws.Range("D24").Value = 252
$ws.Range("E24").Value = -5.952380952380
$ws.Range("F24").Value = 1208
$ws.Range("G24").Value = 1341
$ws.Range("H24").Value = -9.917971662938
$ws.Range("I24").Value = 16160
$ws.Range("J24").Value = 16771
$ws.Range("K24").Value = -3.643193608013
$ws.Range("L24").Value = 34.286189130796
$ws.Range("M24").Value = 38.119658119658
$ws.Range("C25").Value = 201
$ws.Range("D25").Value = 160
$ws.Range("E25").Value = 25.625
$ws.Range("F25").Value = 783
$ws.Range("G25").Value = 706
$ws.Range("H25").Value = 10.906515580736
$ws.Range("I25").Value = 9513
$ws.Range("J25").Value = 9003
$ws.Range("K25").Value = 5.664778407197
$ws.Range("L25").Value = 19.600201156650
$ws.Range("M25").Value = -6.496953017495
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 116.666666666667
$ws.Range("F26").Value = 33
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 6.451612903225
$ws.Range("I26").Value = 571
$ws.Range("J26").Value = 595
$ws.Range("K26").Value = -4.033613445378
$ws.Range("L26").Value = 4.387568555758
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 69
$ws.Range("G27").Value = 57
$ws.Range("H27").Value = 21.052631578947
$ws.Range("I27").Value = 955
$ws.Range("J27").Value = 833
$ws.Range("K27").Value = 14.645858343337
$ws.Range("L27").Value = 12.750885478158
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 13
$ws.Range("E28").Value = -69.230769230769
$ws.Range("F28").Value = 19
$ws.Range("G28").Value = 31
$ws.Range("H28").Value = -38.709677419354
$ws.Range("I28").Value = 355
$ws.Range("J28").Value = 441
$ws.Range("K28").Value = -19.501133786848
$ws.Range("L28").Value = -36.151079136690
$ws.Range("M28").Value = -20.581655480984
$ws.Range("N28").Value = -72.921434019832
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 13
$ws.Range("E29").Value = -69.230769230769
$ws.Range("F29").Value = 18
$ws.Range("G29").Value = 29
$ws.Range("H29").Value = -37.931034482758
$ws.Range("I29").Value = 297
$ws.Range("J29").Value = 376
$ws.Range("K29").Value = -21.010638297872
$ws.Range("L29").Value = -36.538461538461
$ws.Range("M29").Value = -20.8
$ws.Range("N29").Value = -74.915540540540
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 42
$ws.Range("K30").Value = -38.095238095238
